# Balance the quest exp: lower the base XP values for levels 1-3,
# which ripple through the shared formula in column B for subsequent rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp")

$ws.Range("B4").Value = 25
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 40

$excel.Calculate()

# Update the active selection to match the saved view state.
$ws.Range("E6").Select()
